$p = $ppt.ActivePresentation

# Locate the slide / shape that holds the SmartArt diagram ("Normalization
# and clean up of the data ..." staggered process) and fix the wording from
# "clean up" to "cleanup" inside its first node.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasSmartArt) {
            $smartArt = $shape.SmartArt
            for ($ni = 1; $ni -le $smartArt.AllNodes.Count; $ni++) {
                $node = $smartArt.AllNodes.Item($ni)
                if ($node.TextFrame2.TextRange.Text -like "*clean up*") {
                    $node.TextFrame2.TextRange.Text = $node.TextFrame2.TextRange.Text -replace "clean up", "cleanup"
                }
            }
        }
    }
}
